$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the three new sheets in the right tab order:
#    domains, connections, email, toAddresses, backup
# ---------------------------------------------------------------------------
$domainsSheet = $wb.Worksheets.Item("domains")
$connectionsSheet = $wb.Worksheets.Add($null, $domainsSheet)
$connectionsSheet.Name = "connections"

$emailSheet = $wb.Worksheets.Item("email")
$toAddressesSheet = $wb.Worksheets.Add($null, $emailSheet)
$toAddressesSheet.Name = "toAddresses"

$backupSheet = $wb.Worksheets.Add($null, $toAddressesSheet)
$backupSheet.Name = "backup"

# ---------------------------------------------------------------------------
# 2. domains sheet - fill in rows 2 & 3, drop the now-unused row 4
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("domains")

# Give F2:F3 (and later G2:G3) the same look as the rest of the data rows
# (bold green 9pt Courier New, vertically centred) by copying the style that
# A2 already carries instead of rebuilding it font property by font property.
$ws.Range("A2").Copy()
$ws.Range("F2:F3").PasteSpecial(-4122)

$ws.Range("A2").Value = "PKHMapUnits"
$ws.Range("B2").Value = "PKH list of map units"
$ws.Range("C2").Value = "\\sasdfafdsdsf\PKH_LMU_MASTER.xlsx"
$ws.Range("D2").Value = "PKH_LMU$"
$ws.Range("E2").Value = "mapunit"
$ws.Range("F2").Value = "DomainDesc"

$ws.Range("A3").Value = "HPGPS_Purpose"
$ws.Range("B3").Value = "Purpose for HPGPS data"
$ws.Range("C3").Value = "\\Idsfds\AttributeDomains_MASTER.gdb\HPGPS_Purpose"
$ws.Range("E3").Value = "Code"
$ws.Range("F3").Value = "DomainDescrip"

$ws.Rows.Item(4).Delete()

# New "DBs" column values, styled with the non-bold black Courier New look:
# start from the existing bold/green style and only flip the two properties
# that differ, so we don't leave a trail of unused intermediate fonts.
$ws.Range("A2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "[db1,db2,db3]"
$ws.Range("G2").Font.FontStyle = "Regular"
$ws.Range("G2").Font.Color = 0

$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = "[db1]"

$ws.Columns.Item(2).ColumnWidth = 23.5882
$ws.Columns.Item(3).ColumnWidth = 109.5882
$ws.Columns.Item(4).ColumnWidth = 21.5882
$ws.Columns.Item(7).ColumnWidth = 25.5928

$ws.Range("A8").Select()

# ---------------------------------------------------------------------------
# 3. connections sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("connections")

$ws.Range("A1").Value = "dbsNames"
$ws.Range("B1").Value = "sdeConnections"

$domainsSheet.Range("A2").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)

$ws.Range("A2").Value = "db1"
$ws.Range("B2").Value = "Database Connections\Connection to asdfsdfsdfs_SDE.sde"

$ws.Range("A3").Value = "db2"
$ws.Range("B3").Value = "Database Connections\Connection to dfsdfsdfsfdsfdfs_DBO.sde"

$ws.Range("A4").Value = "db3"
$ws.Range("B4").Value = "Database Connections\Connection to asdasddfs_DBO.sde"

$ws.Columns.Item(1).ColumnWidth = 11.5882
$ws.Columns.Item(2).ColumnWidth = 79.5882

$ws.Range("B8").Select()

# ---------------------------------------------------------------------------
# 4. toAddresses sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("toAddresses")

$ws.Range("A1").Value = "email"

$domainsSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "sadasds@gmail.com"
$ws.Range("A2").Font.FontStyle = "Regular"
$ws.Range("A2").Font.Color = 0

$domainsSheet.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = "asdasdsa@usgs.gov"

$ws.Columns.Item(1).ColumnWidth = 17.5882

$ws.Range("M7").Select()

# ---------------------------------------------------------------------------
# 5. backup sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("backup")

$ws.Range("A1").Value = "path"

$domainsSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "\\swdasdas\DomainTableBackups.gdb"

$ws.Range("G6").Select()
$ws.Activate()
